# Applies the commit's changes to "ÁREA GESTIÓN CORPORATIVA" (sheet 1).
# The "CONSOLIDADO ACADÉMICO" sheet only contains formulas that reference
# this sheet (and "ÁREA TECH", which is unchanged), so it recalculates
# automatically once these values/formulas are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 - alumnado ACTIVO
$ws.Range("B3").Formula = "=24+20+33+39+17+8"
$ws.Range("E3").Formula = "=10+15+4+5"

# Row 4 - alumnado INACTIVO
$ws.Range("B4").Formula = "=20+11+1"
$ws.Range("E4").Formula = "=12+6"

# Row 5 - Exito academico
$ws.Range("C5").Formula = "=33/109"
$ws.Range("F5").Formula = "=16/43"

# Row 6 - Absentismo
$ws.Range("C6").Value = 0.0284
$ws.Range("F6").Value = 0.0303

# Row 7 - Riesgo
$ws.Range("C7").Value = 0.01
$ws.Range("F7").Value = 0.02

# Row 10 - Cierre Expediente academico
$ws.Range("C10").Formula = "=17/67"
$ws.Range("F10").Formula = "=10/43"

# Row 11 - Satisfaccion Alumnado
$ws.Range("C11").Value = 0.85
$ws.Range("F11").Value = 0.84

# Row 12 - Resenas
$ws.Range("C12").Formula = "=39/122"

# Row 18 - TOEIC
$ws.Range("C18").Value = 4
$ws.Range("F18").Value = 3

# Row 20 - Bizneo / ISO 27001
$ws.Range("C20").Value = 36
$ws.Range("F20").Value = 7

# Row 21 - CPCC
$ws.Range("F21").Value = 16

# Row 22 - Excel / Cumplen
$ws.Range("C22").Value = 29
$ws.Range("F22").Value = 2

# Row 23 - HRider / Global Suite
$ws.Range("C23").Value = 30
$ws.Range("F23").Value = 42

# Column A width shrank from 8.28515625 to 2.7109375 (OOXML "width" units).
# The COM ColumnWidth property here snaps to an internal 1/6-character grid,
# so 1.8 is the closest input that lands on the nearest achievable width.
$ws.Columns.Item(1).ColumnWidth = 1.8

# Selection moved from D5 to G23.
$ws.Range("G23").Select()
